# The original commit ("updates and data cleaning") widened column A so
# the (long) author lists are fully visible - i.e. an AutoFit/"best fit"
# was applied to column A after the data was cleaned up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
